# Slide 7 ("The New Developer") layout + tagline update.
#
# PowerPoint COM measures Shape.Left/Top/Width/Height in points (a
# `Single`/float32), while the source OOXML stores offsets/extents in
# EMU (914400 EMU per inch, 12700 EMU per point). A naive
# EMU / 12700 division can land one EMU off after the float32 round-trip,
# so convert via a small helper that nudges the point value until it
# reproduces the exact target EMU once PowerPoint stores it as a Single.
function ConvertEmuToPoints($emu) {
    $emuPerPoint = 12700
    $pt = $emu / $emuPerPoint
    for ($i = 0; $i -lt 20000; $i++) {
        $roundTripped = [int64]([single]$pt * $emuPerPoint)
        if ($roundTripped -eq $emu) {
            return $pt
        }
        if ($roundTripped -lt $emu) {
            $pt = $pt + 0.0000001
        } else {
            $pt = $pt - 0.0000001
        }
    }
    return $pt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Text 0 — move down slightly.
$shape1 = $s.Shapes.Item(1)
$shape1.Top = ConvertEmuToPoints 908298

# Text 1 — move down slightly.
$shape2 = $s.Shapes.Item(2)
$shape2.Top = ConvertEmuToPoints 1721048

# Text 2 — move down slightly.
$shape3 = $s.Shapes.Item(3)
$shape3.Top = ConvertEmuToPoints 3008709

# Text 3 — move down, shrink box, and update the tagline copy.
$shape4 = $s.Shapes.Item(4)
$shape4.Top = ConvertEmuToPoints 3719661
$shape4.Width = ConvertEmuToPoints 6263006
$shape4.Height = ConvertEmuToPoints 337691
$shape4.TextFrame.TextRange.Text = "Less time writing code. More time shaping outcomes."
